$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Josephines_Match")

# Fill the match matrix (rows 2-11, columns B-K) with 0 for every
# cell that doesn't already contain the match indicator 1.
for ($r = 2; $r -le 11; $r++) {
    for ($c = 2; $c -le 11; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.Value2 -ne 1) {
            $cell.Value = 0
        }
    }
}

# Update the active selection on the sheet to match the saved state.
$ws.Range("L13").Select()
